# Edit script: applies the two changes captured by the target diff.
#
# 1) Slide 16's table switches from the deck's custom "Table_0" style
#    ({C4D4103C-875B-4F71-BCD6-E43640027DB9}) to the built-in table
#    style {66F87CF1-4007-4F60-8CFD-65003AA3933B}.
#
# 2) The presentation's (slide-side) theme colour palette switches from
#    the "Integral" palette to the stock "Office" palette (this is the
#    substantive, reachable part of the theme1.xml <-> theme2.xml swap
#    recorded in the diff -- theme1.xml is the only theme part exposed
#    for editing through the PowerPoint object model; the Notes
#    Master's independent theme part is not reachable via COM).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 --------------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{66F87CF1-4007-4F60-8CFD-65003AA3933B}")
    }
}

# --- 2) Theme colour palette: Integral -> Office ----------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      #000000
$tcs.Item(2).RGB  = 16777215   # lt1      #FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      #44546A
$tcs.Item(4).RGB  = 15132391   # lt2      #E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  #5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  #ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  #A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  #FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  #4472C4
$tcs.Item(10).RGB = 4697456    # accent6  #70AD47
$tcs.Item(11).RGB = 12673797   # hlink    #0563C1
$tcs.Item(12).RGB = 7491477    # folHlink #954F72
